$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 297.33334
$ws.Range("I28").Value = 297.33334
$ws.Range("K28").Value = 297.33334
$ws.Range("M28").Value = 187.66666
$ws.Range("H39").Value = 1630.0667
$ws.Range("I39").Value = 1077.091
$ws.Range("K39").Value = 3231.273
$ws.Range("M39").Value = -2935.273
$ws.Range("H40").Value = 10612.923
$ws.Range("J40").Value = 3472.25
$ws.Range("L40").Value = 3472.25
$ws.Range("N40").Value = -3822.25
$ws.Range("H42").Value = 401.1
$ws.Range("I42").Value = 241.85715
$ws.Range("J42").Value = 772.6667
$ws.Range("K42").Value = 725.5714499999999
$ws.Range("L42").Value = 2318.0001
$ws.Range("M42").Value = -495.5714499999999
$ws.Range("N42").Value = -2778.0001
$ws.Range("H43").Value = 515616.12
$ws.Range("J43").Value = 823787.8
$ws.Range("L43").Value = 823787.8
$ws.Range("N43").Value = -823925.8
$ws.Range("H55").Value = 308.75
$ws.Range("I55").Value = 342.6
$ws.Range("J55").Value = 252.33333
$ws.Range("K55").Value = 342.6
$ws.Range("L55").Value = 252.33333
$ws.Range("M55").Value = -128.6
$ws.Range("N55").Value = -680.3333299999999
$ws.Range("H58").Value = 4909.3335
$ws.Range("I58").Value = 1046.5
$ws.Range("J58").Value = 7999.6
$ws.Range("K58").Value = 3139.5
$ws.Range("L58").Value = 23998.8
$ws.Range("M58").Value = -2989.5
$ws.Range("N58").Value = -24298.8
$ws.Range("H62").Value = 27999.15
$ws.Range("I62").Value = 1770.5
$ws.Range("J62").Value = 54227.8
$ws.Range("K62").Value = 1770.5
$ws.Range("L62").Value = 54227.8
$ws.Range("M62").Value = -1146.5
$ws.Range("N62").Value = -55475.8
$ws.Range("H65").Value = 27999.15
$ws.Range("I65").Value = 1770.5
$ws.Range("J65").Value = 54227.8
$ws.Range("K65").Value = 8852.5
$ws.Range("L65").Value = 271139
$ws.Range("M65").Value = -5732.5
$ws.Range("N65").Value = -277379
$ws.Range("H70").Value = 1473.5714
$ws.Range("I70").Value = 1218.8
$ws.Range("J70").Value = 1553.1875
$ws.Range("K70").Value = 3656.4
$ws.Range("L70").Value = 4659.5625
$ws.Range("M70").Value = -3386.4
$ws.Range("N70").Value = -5199.5625
$ws.Range("H73").Value = 1473.5714
$ws.Range("I73").Value = 1218.8
$ws.Range("J73").Value = 1553.1875
$ws.Range("K73").Value = 3656.4
$ws.Range("L73").Value = 4659.5625
$ws.Range("M73").Value = -2720.4
$ws.Range("N73").Value = -6531.5625
$ws.Range("H74").Value = 35719744
$ws.Range("I74").Value = 71431290
$ws.Range("J74").Value = 8201.143
$ws.Range("K74").Value = 71431290
$ws.Range("L74").Value = 8201.143
$ws.Range("M74").Value = -71430354
$ws.Range("N74").Value = -10073.143
$ws.Range("H77").Value = 35719744
$ws.Range("I77").Value = 71431290
$ws.Range("J77").Value = 8201.143
$ws.Range("K77").Value = 357156450
$ws.Range("L77").Value = 41005.715
$ws.Range("M77").Value = -357151770
$ws.Range("N77").Value = -50365.715
$ws.Range("H80").Value = 56938.777
$ws.Range("J80").Value = 58849.855
$ws.Range("L80").Value = 176549.565
$ws.Range("N80").Value = -178545.565
$ws.Range("H83").Value = 56938.777
$ws.Range("J83").Value = 58849.855
$ws.Range("L83").Value = 529648.6950000001
$ws.Range("N83").Value = -539632.6950000001
$ws.Range("H86").Value = 66139210
$ws.Range("I86").Value = 83334940
$ws.Range("K86").Value = 83334940
$ws.Range("M86").Value = -83333817
$ws.Range("H87").Value = 59995
$ws.Range("J87").Value = 59995
$ws.Range("L87").Value = 59995
$ws.Range("N87").Value = -62491
$ws.Range("H89").Value = 66139210
$ws.Range("I89").Value = 83334940
$ws.Range("K89").Value = 416674700
$ws.Range("M89").Value = -416669084
$ws.Range("H90").Value = 59995
$ws.Range("J90").Value = 59995
$ws.Range("L90").Value = 179985
$ws.Range("N90").Value = -192465
$ws.Range("H98").Value = 9185.799999999999
$ws.Range("I98").Value = 9268.23
$ws.Range("J98").Value = 8650
$ws.Range("K98").Value = 9268.23
$ws.Range("L98").Value = 8650
$ws.Range("M98").Value = -7770.23
$ws.Range("N98").Value = -11646
$ws.Range("H115").Value = 955.55554
$ws.Range("I115").Value = 1026.25
$ws.Range("K115").Value = 3078.75
$ws.Range("M115").Value = -1511.75
$ws.Range("H122").Value = 9185.799999999999
$ws.Range("I122").Value = 9268.23
$ws.Range("J122").Value = 8650
$ws.Range("K122").Value = 27804.69
$ws.Range("L122").Value = 25950
$ws.Range("M122").Value = -25354.69
$ws.Range("N122").Value = -30850
$ws.Range("H132").Value = 2564.8386
$ws.Range("I132").Value = 2564.8386
$ws.Range("K132").Value = 7694.5158
$ws.Range("M132").Value = -5164.5158
$ws.Range("H135").Value = 910658.5600000001
$ws.Range("I135").Value = 1429427.8
$ws.Range("J135").Value = 2812.5
$ws.Range("K135").Value = 12864850.2
$ws.Range("L135").Value = 25312.5
$ws.Range("M135").Value = -12862315.2
$ws.Range("N135").Value = -30382.5
$ws.Range("H138").Value = 5520.6777
$ws.Range("J138").Value = 6922.878
$ws.Range("L138").Value = 20768.634
$ws.Range("N138").Value = -31048.634
$ws.Range("H141").Value = 2296.6667
$ws.Range("I141").Value = 2296.6667
$ws.Range("K141").Value = 6890.000100000001
$ws.Range("M141").Value = -1710.000100000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 4704.6665
$ws.Range("J14").Value = 4704.6665
$ws.Range("L14").Value = 4704.6665
$ws.Range("N14").Value = -5054.6665
$ws.Range("H16").Value = 2258
$ws.Range("J16").Value = 4000
$ws.Range("L16").Value = 4000
$ws.Range("N16").Value = -4574
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996
$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984
$ws.Range("H97").Value = 2981049.8
$ws.Range("I97").Value = 608.7368
$ws.Range("K97").Value = 608.7368
$ws.Range("M97").Value = -112.7368
$ws.Range("H122").Value = 2672.9148
$ws.Range("I122").Value = 2112.6943
$ws.Range("K122").Value = 6338.0829
$ws.Range("M122").Value = -3888.0829
$ws.Range("H132").Value = 5476.204
$ws.Range("I132").Value = 4235.303
$ws.Range("J132").Value = 8035.5625
$ws.Range("K132").Value = 12705.909
$ws.Range("L132").Value = 24106.6875
$ws.Range("M132").Value = -10175.909
$ws.Range("N132").Value = -29166.6875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 48754.5
$ws.Range("J28").Value = 48754.5
$ws.Range("L28").Value = 48754.5
$ws.Range("N28").Value = -49342.5
$ws.Range("H36").Value = 425
$ws.Range("I36").Value = 425
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 425
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 109
$ws.Range("H64").Value = 22223376
$ws.Range("I64").Value = 41667520
$ws.Range("J64").Value = 1495
$ws.Range("K64").Value = 41667520
$ws.Range("L64").Value = 1495
$ws.Range("M64").Value = -41667295
$ws.Range("N64").Value = -1945
$ws.Range("H67").Value = 22223376
$ws.Range("I67").Value = 41667520
$ws.Range("J67").Value = 1495
$ws.Range("K67").Value = 41667520
$ws.Range("L67").Value = 1495
$ws.Range("M67").Value = -41666740
$ws.Range("N67").Value = -3055
$ws.Range("H82").Value = 4723.25
$ws.Range("I82").Value = 4723.25
$ws.Range("K82").Value = 4723.25
$ws.Range("M82").Value = -4340.25
$ws.Range("H85").Value = 4723.25
$ws.Range("I85").Value = 4723.25
$ws.Range("K85").Value = 4723.25
$ws.Range("M85").Value = -3397.25
$ws.Range("H86").Value = 53296.5
$ws.Range("I86").Value = 85610.25
$ws.Range("J86").Value = 4825.875
$ws.Range("K86").Value = 85610.25
$ws.Range("L86").Value = 4825.875
$ws.Range("M86").Value = -84487.25
$ws.Range("N86").Value = -7071.875
$ws.Range("H89").Value = 53296.5
$ws.Range("I89").Value = 85610.25
$ws.Range("J89").Value = 4825.875
$ws.Range("K89").Value = 428051.25
$ws.Range("L89").Value = 24129.375
$ws.Range("M89").Value = -422435.25
$ws.Range("N89").Value = -35361.375
$ws.Range("H94").Value = 331.03333
$ws.Range("I94").Value = 197.14285
$ws.Range("J94").Value = 643.44446
$ws.Range("K94").Value = 197.14285
$ws.Range("L94").Value = 643.44446
$ws.Range("M94").Value = 253.85715
$ws.Range("N94").Value = -1545.44446
$ws.Range("H105").Value = 5034.1113
$ws.Range("I105").Value = 2450.5
$ws.Range("K105").Value = 2450.5
$ws.Range("M105").Value = -703.5
$ws.Range("H134").Value = 7157
$ws.Range("I134").Value = 3743.1333
$ws.Range("J134").Value = 10169.235
$ws.Range("K134").Value = 11229.3999
$ws.Range("L134").Value = 30507.705
$ws.Range("M134").Value = -8694.3999
$ws.Range("N134").Value = -35577.705
$ws.Range("N36").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 45454736
$ws.Range("I7").Value = 48.444443
$ws.Range("K7").Value = 48.444443
$ws.Range("M7").Value = 64.55555699999999
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 466.66666
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 466.66666
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -116.66666
$ws.Range("N22").Value = -1700
$ws.Range("H31").Value = 10068.296
$ws.Range("I31").Value = 4501.077
$ws.Range("K31").Value = 4501.077
$ws.Range("M31").Value = -4206.077
$ws.Range("H34").Value = 10068.296
$ws.Range("I34").Value = 4501.077
$ws.Range("K34").Value = 4501.077
$ws.Range("M34").Value = -4299.077
$ws.Range("H58").Value = 7478.5586
$ws.Range("I58").Value = 3253.1428
$ws.Range("J58").Value = 8574.037
$ws.Range("K58").Value = 3253.1428
$ws.Range("L58").Value = 8574.037
$ws.Range("M58").Value = -3050.1428
$ws.Range("N58").Value = -8980.037
$ws.Range("H62").Value = 20837096
$ws.Range("I62").Value = 25002516
$ws.Range("K62").Value = 25002516
$ws.Range("M62").Value = -25001892
$ws.Range("H65").Value = 20837096
$ws.Range("I65").Value = 25002516
$ws.Range("K65").Value = 125012580
$ws.Range("M65").Value = -125009460
$ws.Range("H86").Value = 30313528
$ws.Range("I86").Value = 10113307
$ws.Range("K86").Value = 10113307
$ws.Range("M86").Value = -10112184
$ws.Range("H89").Value = 30313528
$ws.Range("I89").Value = 10113307
$ws.Range("K89").Value = 50566535
$ws.Range("M89").Value = -50560919
$ws.Range("H99").Value = 4634.6
$ws.Range("I99").Value = 3200
$ws.Range("K99").Value = 3200
$ws.Range("M99").Value = -1702
$ws.Range("H105").Value = 4465790
$ws.Range("I105").Value = 5953150.5
$ws.Range("K105").Value = 5953150.5
$ws.Range("M105").Value = -5951403.5
$ws.Range("H126").Value = 4634.6
$ws.Range("I126").Value = 3200
$ws.Range("K126").Value = 9600
$ws.Range("M126").Value = -7130
$ws.Range("H132").Value = 5840.1797
$ws.Range("I132").Value = 3799.48
$ws.Range("K132").Value = 11398.44
$ws.Range("M132").Value = -8868.440000000001
$ws.Range("H134").Value = 6462.6855
$ws.Range("I134").Value = 2954.6428
$ws.Range("K134").Value = 8863.928400000001
$ws.Range("M134").Value = -6328.928400000001
$ws.Range("H136").Value = 7478.5586
$ws.Range("I136").Value = 3253.1428
$ws.Range("J136").Value = 8574.037
$ws.Range("K136").Value = 9759.428400000001
$ws.Range("L136").Value = 25722.111
$ws.Range("M136").Value = -7209.428400000001
$ws.Range("N136").Value = -30822.111

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1715.1
$ws.Range("I25").Value = 450.2
$ws.Range("J25").Value = 2980
$ws.Range("K25").Value = 1350.6
$ws.Range("L25").Value = 8940
$ws.Range("M25").Value = -1181.6
$ws.Range("N25").Value = -9278
$ws.Range("H30").Value = 1715.1
$ws.Range("I30").Value = 450.2
$ws.Range("J30").Value = 2980
$ws.Range("K30").Value = 1350.6
$ws.Range("L30").Value = 8940
$ws.Range("M30").Value = -1248.6
$ws.Range("N30").Value = -9144
$ws.Range("H34").Value = 4651.647
$ws.Range("J34").Value = 5633.357
$ws.Range("L34").Value = 16900.071
$ws.Range("N34").Value = -17068.071
$ws.Range("H38").Value = 47.7
$ws.Range("I38").Value = 43.75
$ws.Range("J38").Value = 50.333332
$ws.Range("K38").Value = 131.25
$ws.Range("L38").Value = 150.999996
$ws.Range("M38").Value = 215.75
$ws.Range("N38").Value = -844.999996
$ws.Range("H63").Value = 1914
$ws.Range("J63").Value = 1914
$ws.Range("L63").Value = 5742
$ws.Range("N63").Value = -7240
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("H66").Value = 1914
$ws.Range("J66").Value = 1914
$ws.Range("L66").Value = 17226
$ws.Range("N66").Value = -24714
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("H68").Value = 14774
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 21161
$ws.Range("K68").Value = 6000
$ws.Range("L68").Value = 63483
$ws.Range("M68").Value = -5189
$ws.Range("N68").Value = -65105
$ws.Range("H71").Value = 14774
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 21161
$ws.Range("K71").Value = 18000
$ws.Range("L71").Value = 190449
$ws.Range("M71").Value = -13944
$ws.Range("N71").Value = -198561
$ws.Range("H75").Value = 95248440
$ws.Range("I75").Value = 83336440
$ws.Range("J75").Value = 111131110
$ws.Range("K75").Value = 250009320
$ws.Range("L75").Value = 333393330
$ws.Range("M75").Value = -250008322
$ws.Range("N75").Value = -333395326
$ws.Range("H78").Value = 95248440
$ws.Range("I78").Value = 83336440
$ws.Range("J78").Value = 111131110
$ws.Range("K78").Value = 750027960
$ws.Range("L78").Value = 1000179990
$ws.Range("M78").Value = -750022968
$ws.Range("N78").Value = -1000189974
$ws.Range("H82").Value = 9982.200000000001
$ws.Range("I82").Value = 7752.75
$ws.Range("K82").Value = 23258.25
$ws.Range("M82").Value = -22852.25
$ws.Range("H85").Value = 9982.200000000001
$ws.Range("I85").Value = 7752.75
$ws.Range("K85").Value = 23258.25
$ws.Range("M85").Value = -21854.25
$ws.Range("H98").Value = 55555772
$ws.Range("I98").Value = 299.5
$ws.Range("J98").Value = 100000150
$ws.Range("K98").Value = 898.5
$ws.Range("L98").Value = 300000450
$ws.Range("M98").Value = 599.5
$ws.Range("N98").Value = -300003446
$ws.Range("H107").Value = 711
$ws.Range("J107").Value = 728.6667
$ws.Range("L107").Value = 2186.0001
$ws.Range("N107").Value = -6026.0001
$ws.Range("H108").Value = 100
$ws.Range("I108").Value = 100
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 300
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 2580
$ws.Range("H113").Value = 3026.1785
$ws.Range("I113").Value = 823.2222
$ws.Range("K113").Value = 2469.6666
$ws.Range("M113").Value = -299.6666
$ws.Range("H133").Value = 2999.6667
$ws.Range("I133").Value = 1800
$ws.Range("J133").Value = 3599.5
$ws.Range("K133").Value = 5400
$ws.Range("L133").Value = 10798.5
$ws.Range("M133").Value = -340
$ws.Range("N133").Value = -20918.5
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("N108").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 20001
$ws.Range("J48").Value = 20001
$ws.Range("L48").Value = 20001
$ws.Range("N48").Value = -20971
$ws.Range("H70").Value = 11999.857
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 11999.857
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 11999.857
$ws.Range("N70").Value = -12539.857
$ws.Range("H73").Value = 11999.857
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 11999.857
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 11999.857
$ws.Range("N73").Value = -13871.857
$ws.Range("H97").Value = 1159.6666
$ws.Range("I97").Value = 1021.931
$ws.Range("J97").Value = 1466.9231
$ws.Range("K97").Value = 1021.931
$ws.Range("L97").Value = 1466.9231
$ws.Range("M97").Value = -525.931
$ws.Range("N97").Value = -2458.9231
$ws.Range("H102").Value = 2010.8572
$ws.Range("I102").Value = 2023.9
$ws.Range("K102").Value = 2023.9
$ws.Range("M102").Value = -401.9000000000001
$ws.Range("H113").Value = 4537.0527
$ws.Range("I113").Value = 2427.7273
$ws.Range("K113").Value = 2427.7273
$ws.Range("M113").Value = -257.7273
$ws.Range("H126").Value = 45457036
$ws.Range("I126").Value = 166669140
$ws.Range("J126").Value = 2499.5
$ws.Range("K126").Value = 500007420
$ws.Range("L126").Value = 7498.5
$ws.Range("M126").Value = -500004950
$ws.Range("N126").Value = -12438.5
$ws.Range("H132").Value = 4820.143
$ws.Range("I132").Value = 1456.4286
$ws.Range("J132").Value = 11547.571
$ws.Range("K132").Value = 4369.2858
$ws.Range("L132").Value = 34642.713
$ws.Range("M132").Value = -1839.2858
$ws.Range("N132").Value = -39702.713
$ws.Range("H140").Value = 25780
$ws.Range("J140").Value = 25780
$ws.Range("L140").Value = 25780
$ws.Range("N140").Value = -36140
$ws.Range("H141").Value = 9993
$ws.Range("I141").Value = 9993
$ws.Range("K141").Value = 9993
$ws.Range("M141").Value = -4813
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3721
$ws.Range("I7").Value = 2819.8948
$ws.Range("K7").Value = 2819.8948
$ws.Range("M7").Value = -2707.8948
$ws.Range("H55").Value = 492.33334
$ws.Range("I55").Value = 105.5
$ws.Range("K55").Value = 105.5
$ws.Range("M55").Value = 67.5
$ws.Range("H126").Value = 3721
$ws.Range("I126").Value = 2819.8948
$ws.Range("K126").Value = 8459.6844
$ws.Range("M126").Value = -5989.6844
$ws.Range("H132").Value = 6598.244
$ws.Range("I132").Value = 3529.25
$ws.Range("J132").Value = 10930.941
$ws.Range("K132").Value = 10587.75
$ws.Range("L132").Value = 32792.823
$ws.Range("M132").Value = -8057.75
$ws.Range("N132").Value = -37852.823
$ws.Range("H136").Value = 9859.516
$ws.Range("I136").Value = 4587.1113
$ws.Range("K136").Value = 13761.3339
$ws.Range("M136").Value = -11211.3339

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 12000
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("H62").Value = 8640
$ws.Range("I62").Value = 7652.4
$ws.Range("K62").Value = 7652.4
$ws.Range("M62").Value = -7028.4
$ws.Range("H65").Value = 8640
$ws.Range("I65").Value = 7652.4
$ws.Range("K65").Value = 38262
$ws.Range("M65").Value = -35142
$ws.Range("H100").Value = 425
$ws.Range("I100").Value = 425
$ws.Range("K100").Value = 850
$ws.Range("M100").Value = -309
$ws.Range("H122").Value = 12925095
$ws.Range("I122").Value = 18667844
$ws.Range("J122").Value = 3908.9167
$ws.Range("K122").Value = 56003532
$ws.Range("L122").Value = 11726.7501
$ws.Range("M122").Value = -56001082
$ws.Range("N122").Value = -16626.7501
$ws.Range("H132").Value = 35705.055
$ws.Range("I132").Value = 8918.1875
$ws.Range("K132").Value = 26754.5625
$ws.Range("M132").Value = -24224.5625
$ws.Range("N47").ClearContents()
